# The "LoadTime" sheet previously had a header row ("Page Load") followed by
# three labels ("Trucks", "Dashboard", "Login") spread across a wide,
# otherwise-empty A1:J4 block. The new export tightens this down to a
# compact two-column report: column A keeps the labels, column B carries
# the actual loading-time figures (a date stamp for the header, and the
# numeric timings for each row), and the now-unused C:J block is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused C:J columns for rows 1-4 (shrinks the used range down
# to A1:B4).
$ws.Range("C1:J4").ClearContents() | Out-Null

# Populate the new column B with the loading-time data.
$ws.Range("B1").Value = " Sep 16"
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 9
$ws.Range("B4").Value = 30

# Match the author's final selection/cursor position.
$ws.Range("B8").Select() | Out-Null
